$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "GANCI CORRADO" is replaced by a new technician "GIONFRIDDO ANDREA"
# (the stock count in C3 is unchanged)
$ws.Range("B3").Value = "GIONFRIDDO ANDREA"

# Append three new stock-count rows (8, 9, 10), copying the formatting
# (date style + font style) from the last existing data row (7)
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C10").PasteSpecial(-4122)

$ws.Range("A8").Value = 45887
$ws.Range("B8").Value = "MANCARELLA SALVATORE"
$ws.Range("C8").Value = 5

$ws.Range("A9").Value = 45887
$ws.Range("B9").Value = "AVOLA IVAN"
$ws.Range("C9").Value = 1

$ws.Range("A10").Value = 45887
$ws.Range("B10").Value = "GOLINO KEVIN"
$ws.Range("C10").Value = 15

# Column B widens slightly to fit the longest new name ("MANCARELLA SALVATORE")
$ws.Columns("B").ColumnWidth = 23.3

# Move the active selection to B12, matching the author's last cursor position
[void]$ws.Range("B12").Select()
